$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.821.37'
$ws.Range('E2').Value = '  -0.12%  '
$ws.Range('D3').Value = '2.078.21'
$ws.Range('E3').Value = '  -1.27%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '233.40'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.32%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.626'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.17%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '59.09'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +2.11%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.394'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +0.84%  '
$ws.Range('E10').Value = '  +0.53%  '
$ws.Range('E11').Value = '  +1.15%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '14.81'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +1.16%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '21.17'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -0.99%  '
$ws.Range('E14').Value = '  -0.61%  '
$ws.Range('E15').Value = '  +1.59%  '
$ws.Range('D16').Value = '2.055.84'
$ws.Range('E16').Value = '  -2.09%  '
$ws.Range('D17').Value = '37.723.94'
$ws.Range('E17').Value = '  -0.25%  '
$ws.Range('E18').Value = '  -1.13%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '71.62'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +1.04%  '
$ws.Range('D20').Value = '0.0₃0840'
$ws.Range('E20').Value = '  +1.68%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '228.41'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.17%  '
$ws.Range('E22').Value = '  -0.11%  '
$ws.Range('E23').Value = '  +0.44%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.35'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -2.60%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '171.07'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +1.78%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.18'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +2.01%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.136'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -3.15%  '
$ws.Range('E28').Value = '  -0.19%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '19.49'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.25%  '
$ws.Range('E30').Value = '  +1.72%  '
$ws.Range('E31').Value = '  +0.21%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.77'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +2.31%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0635'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.61%  '
$ws.Range('E34').Value = '  -3.69%  '
$ws.Range('E35').Value = '  +0.07%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.41'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -1.36%  '
$ws.Range('E37').Value = '  +0.17%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '5.38'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.89%  '
$ws.Range('E39').Value = '  -1.03%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '99.50'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +1.70%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0216'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.08%  '
$ws.Range('B42').Value = 'HuobiToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.88'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -2.46%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '16.70'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +5.91%  '
$ws.Range('D44').Value = '1.442.01'
$ws.Range('E44').Value = '  -1.11%  '
$ws.Range('E45').Value = '  -0.88%  '
$ws.Range('E46').Value = '  +2.05%  '
$ws.Range('E47').Value = '  -0.39%  '
$ws.Range('E48').Value = '  +1.38%  '
$ws.Range('E49').Value = '  -1.59%  '
$ws.Range('D50').Value = '2.270.65'
$ws.Range('E50').Value = '  -1.11%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '46.74'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.51%  '
